$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") serial date value updated from 45496 to 45497
# for every data row (rows 2 through 28).
$ws.Range("C2:C28").Value = 45497
